$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the newly computed "model fit" rows (31, 40, 41) ---
# Row 31
$ws.Range("D31").Value = 130685.02
$ws.Range("E31").Value = 130688.02
$ws.Range("F31").Formula = "=IF(D31<E31,""non-pw"",""pw"")"

# Row 40
$ws.Range("D40").Value = 135122.01
$ws.Range("E40").Value = 135126.31
$ws.Range("F40").Formula = "=IF(D40<E40,""non-pw"",""pw"")"

# Row 41
$ws.Range("D41").Value = 142448.33
$ws.Range("E41").Value = 142452.98
$ws.Range("F41").Formula = "=IF(D41<E41,""non-pw"",""pw"")"

# --- Update sheet view: scroll back to top and select E8 ---
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("E8").Select()
